$wb = $excel.ActiveWorkbook

# --- Config sheet: add new row (B6) describing the config server dependency ---
$configSheet = $wb.Worksheets.Item("Config")
$configSheet.Range("B6").Value = "Just add the dependency in the pom.xml for referring a config server"
$configSheet.Range("B6").WrapText = $true

# Make Config the active/selected sheet and put selection on the new cell
$configSheet.Activate()
$configSheet.Range("B6").Select()

# --- Microservices Frameworks sheet: shrink row 12 height (45 -> 30) ---
$frameworksSheet = $wb.Worksheets.Item("Microservices Frameworks")
$frameworksSheet.Rows.Item(12).RowHeight = 30
